$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AA4").Value = 0.13
$ws.Range("AB4").Value = 0.361
$ws.Range("AC4").Value = 0.117
$ws.Range("AE4").Value = 0.082
$ws.Range("AF4").Value = 0.713
$ws.Range("AH4").Value = 0.325
$ws.Range("AI4").Value = 0.656
$ws.Range("AJ4").Value = 0.165
$ws.Range("AK4").Value = 0.406
$ws.Range("AL4").Value = 0.671
$ws.Range("AO4").Value = 0.68
$ws.Range("B4").Value = 0.289
$ws.Range("E4").Value = 0.168
$ws.Range("F4").Value = 0.028
$ws.Range("G4").Value = 0.168
$ws.Range("H4").Value = 0.2
$ws.Range("I4").Value = 0.026
$ws.Range("J4").Value = 0.161
$ws.Range("K4").Value = 0.307
$ws.Range("L4").Value = 0.099
$ws.Range("M4").Value = 0.314
$ws.Range("N4").Value = 0.254
$ws.Range("P4").Value = 0.149
$ws.Range("Q4").Value = 0.478
$ws.Range("R4").Value = 0.222
$ws.Range("T4").Value = 0.244
$ws.Range("V4").Value = 0.289
$ws.Range("W4").Value = 0.242
$ws.Range("X4").Value = 0.043
$ws.Range("Z4").Value = 0.431
$ws.Range("AA5").Value = 0.16
$ws.Range("AB5").Value = 0.4
$ws.Range("AC5").Value = 0.711
$ws.Range("AD5").Value = 0.205
$ws.Range("AE5").Value = 0.453
$ws.Range("AF5").Value = 0.956
$ws.Range("AG5").Value = 0.042
$ws.Range("AH5").Value = 0.206
$ws.Range("AI5").Value = 0.778
$ws.Range("AJ5").Value = 0.173
$ws.Range("AK5").Value = 0.416
$ws.Range("AL5").Value = 0.911
$ws.Range("AM5").Value = 0.081
$ws.Range("AN5").Value = 0.285
$ws.Range("AO5").Value = 0.882
$ws.Range("B5").Value = 0.822
$ws.Range("C5").Value = 0.146
$ws.Range("D5").Value = 0.382
$ws.Range("E5").Value = 0.667
$ws.Range("F5").Value = 0.222
$ws.Range("G5").Value = 0.471
$ws.Range("H5").Value = 0.8
$ws.Range("I5").Value = 0.16
$ws.Range("J5").Value = 0.4
$ws.Range("K5").Value = 0.6
$ws.Range("L5").Value = 0.24
$ws.Range("M5").Value = 0.49
$ws.Range("N5").Value = 0.8
$ws.Range("O5").Value = 0.16
$ws.Range("P5").Value = 0.4
$ws.Range("Q5").Value = 0.533
$ws.Range("R5").Value = 0.249
$ws.Range("S5").Value = 0.499
$ws.Range("T5").Value = 0.511
$ws.Range("U5").Value = 0.25
$ws.Range("V5").Value = 0.5
$ws.Range("W5").Value = 0.733
$ws.Range("X5").Value = 0.196
$ws.Range("Y5").Value = 0.442
$ws.Range("Z5").Value = 0.8
$ws.Range("AC6").Value = 0.201
$ws.Range("AF6").Value = 0.8169999999999999
$ws.Range("AI6").Value = 0.712
$ws.Range("AL6").Value = 0.773
$ws.Range("AO6").Value = 0.767
$ws.Range("B6").Value = 0.428
$ws.Range("E6").Value = 0.268
$ws.Range("H6").Value = 0.32
$ws.Range("K6").Value = 0.406
$ws.Range("N6").Value = 0.386
$ws.Range("Q6").Value = 0.504
$ws.Range("T6").Value = 0.33
$ws.Range("W6").Value = 0.364
$ws.Range("Z6").Value = 0.5600000000000001
$ws.Range("AC7").Value = 0.353
$ws.Range("AF7").Value = 0.895
$ws.Range("AI7").Value = 0.75
$ws.Range("AL7").Value = 0.85
$ws.Range("AO7").Value = 0.832
$ws.Range("B7").Value = 0.601
$ws.Range("E7").Value = 0.418
$ws.Range("H7").Value = 0.5
$ws.Range("K7").Value = 0.504
$ws.Range("N7").Value = 0.5590000000000001
$ws.Range("Q7").Value = 0.521
$ws.Range("T7").Value = 0.419
$ws.Range("W7").Value = 0.521
$ws.Range("Z7").Value = 0.6830000000000001
$ws.Range("AA8").Value = 0.157
$ws.Range("AB8").Value = 0.396
$ws.Range("AC8").Value = 0.596
$ws.Range("AD8").Value = 0.189
$ws.Range("AE8").Value = 0.435
$ws.Range("AF8").Value = 0.879
$ws.Range("AG8").Value = 0.06
$ws.Range("AH8").Value = 0.244
$ws.Range("AI8").Value = 0.77
$ws.Range("AJ8").Value = 0.172
$ws.Range("AK8").Value = 0.415
$ws.Range("AL8").Value = 0.878
$ws.Range("AM8").Value = 0.08599999999999999
$ws.Range("AN8").Value = 0.294
$ws.Range("AO8").Value = 0.842
$ws.Range("B8").Value = 0.742
$ws.Range("C8").Value = 0.149
$ws.Range("D8").Value = 0.386
$ws.Range("E8").Value = 0.5629999999999999
$ws.Range("H8").Value = 0.697
$ws.Range("I8").Value = 0.158
$ws.Range("J8").Value = 0.398
$ws.Range("K8").Value = 0.531
$ws.Range("M8").Value = 0.46
$ws.Range("N8").Value = 0.713
$ws.Range("O8").Value = 0.157
$ws.Range("P8").Value = 0.396
$ws.Range("Q8").Value = 0.509
$ws.Range("S8").Value = 0.484
$ws.Range("T8").Value = 0.445
$ws.Range("W8").Value = 0.662
$ws.Range("X8").Value = 0.182
$ws.Range("Y8").Value = 0.426
$ws.Range("Z8").Value = 0.737
$ws.Range("AA9").Value = 0.229
$ws.Range("AB9").Value = 0.479
$ws.Range("AC9").Value = 0.489
$ws.Range("AD9").Value = 0.25
$ws.Range("AE9").Value = 0.5
$ws.Range("AF9").Value = 0.756
$ws.Range("AG9").Value = 0.185
$ws.Range("AH9").Value = 0.43
$ws.Range("AI9").Value = 0.756
$ws.Range("AJ9").Value = 0.185
$ws.Range("AK9").Value = 0.43
$ws.Range("AL9").Value = 0.822
$ws.Range("AM9").Value = 0.146
$ws.Range("AN9").Value = 0.382
$ws.Range("AO9").Value = 0.778
$ws.Range("B9").Value = 0.644
$ws.Range("C9").Value = 0.229
$ws.Range("D9").Value = 0.479
$ws.Range("E9").Value = 0.444
$ws.Range("F9").Value = 0.247
$ws.Range("G9").Value = 0.497
$ws.Range("H9").Value = 0.578
$ws.Range("I9").Value = 0.244
$ws.Range("J9").Value = 0.494
$ws.Range("K9").Value = 0.444
$ws.Range("L9").Value = 0.247
$ws.Range("M9").Value = 0.497
$ws.Range("N9").Value = 0.6
$ws.Range("O9").Value = 0.24
$ws.Range("P9").Value = 0.49
$ws.Range("Q9").Value = 0.467
$ws.Range("T9").Value = 0.356
$ws.Range("U9").Value = 0.229
$ws.Range("V9").Value = 0.479
$ws.Range("W9").Value = 0.556
$ws.Range("X9").Value = 0.247
$ws.Range("Y9").Value = 0.497
$ws.Range("Z9").Value = 0.644
$ws.Range("AA10").Value = 0.16
$ws.Range("AB10").Value = 0.4
$ws.Range("AC10").Value = 0.6
$ws.Range("AD10").Value = 0.24
$ws.Range("AE10").Value = 0.49
$ws.Range("AF10").Value = 0.956
$ws.Range("AG10").Value = 0.042
$ws.Range("AH10").Value = 0.206
$ws.Range("AI10").Value = 0.778
$ws.Range("AJ10").Value = 0.173
$ws.Range("AK10").Value = 0.416
$ws.Range("AL10").Value = 0.911
$ws.Range("AM10").Value = 0.081
$ws.Range("AN10").Value = 0.285
$ws.Range("AO10").Value = 0.882
$ws.Range("B10").Value = 0.778
$ws.Range("C10").Value = 0.173
$ws.Range("D10").Value = 0.416
$ws.Range("E10").Value = 0.6
$ws.Range("F10").Value = 0.24
$ws.Range("G10").Value = 0.49
$ws.Range("H10").Value = 0.733
$ws.Range("I10").Value = 0.196
$ws.Range("J10").Value = 0.442
$ws.Range("K10").Value = 0.6
$ws.Range("L10").Value = 0.24
$ws.Range("M10").Value = 0.49
$ws.Range("N10").Value = 0.778
$ws.Range("O10").Value = 0.173
$ws.Range("P10").Value = 0.416
$ws.Range("Q10").Value = 0.533
$ws.Range("R10").Value = 0.249
$ws.Range("S10").Value = 0.499
$ws.Range("T10").Value = 0.511
$ws.Range("U10").Value = 0.25
$ws.Range("V10").Value = 0.5
$ws.Range("W10").Value = 0.733
$ws.Range("X10").Value = 0.196
$ws.Range("Y10").Value = 0.442
$ws.Range("Z10").Value = 0.8
$ws.Range("AA11").Value = 0.16
$ws.Range("AB11").Value = 0.4
$ws.Range("AC11").Value = 0.644
$ws.Range("AD11").Value = 0.229
$ws.Range("AE11").Value = 0.479
$ws.Range("AF11").Value = 0.956
$ws.Range("AG11").Value = 0.042
$ws.Range("AH11").Value = 0.206
$ws.Range("AI11").Value = 0.778
$ws.Range("AJ11").Value = 0.173
$ws.Range("AK11").Value = 0.416
$ws.Range("AL11").Value = 0.911
$ws.Range("AM11").Value = 0.081
$ws.Range("AN11").Value = 0.285
$ws.Range("AO11").Value = 0.882
$ws.Range("B11").Value = 0.822
$ws.Range("C11").Value = 0.146
$ws.Range("D11").Value = 0.382
$ws.Range("E11").Value = 0.667
$ws.Range("F11").Value = 0.222
$ws.Range("G11").Value = 0.471
$ws.Range("H11").Value = 0.8
$ws.Range("I11").Value = 0.16
$ws.Range("J11").Value = 0.4
$ws.Range("K11").Value = 0.6
$ws.Range("L11").Value = 0.24
$ws.Range("M11").Value = 0.49
$ws.Range("N11").Value = 0.8
$ws.Range("O11").Value = 0.16
$ws.Range("P11").Value = 0.4
$ws.Range("Q11").Value = 0.533
$ws.Range("R11").Value = 0.249
$ws.Range("S11").Value = 0.499
$ws.Range("T11").Value = 0.511
$ws.Range("U11").Value = 0.25
$ws.Range("V11").Value = 0.5
$ws.Range("W11").Value = 0.733
$ws.Range("X11").Value = 0.196
$ws.Range("Y11").Value = 0.442
$ws.Range("Z11").Value = 0.8
$ws.Range("AA12").Value = 0.299
$ws.Range("AB12").Value = 0.546
$ws.Range("AC12").Value = 2
$ws.Range("AD12").Value = 3.812
$ws.Range("AE12").Value = 1.953
$ws.Range("AF12").Value = 1.233
$ws.Range("AG12").Value = 0.225
$ws.Range("AH12").Value = 0.474
$ws.Range("AJ12").Value = 0.028
$ws.Range("AK12").Value = 0.167
$ws.Range("AL12").Value = 1.098
$ws.Range("AM12").Value = 0.08799999999999999
$ws.Range("AN12").Value = 0.297
$ws.Range("AO12").Value = 1.12
$ws.Range("B12").Value = 1.378
$ws.Range("C12").Value = 0.668
$ws.Range("D12").Value = 0.8169999999999999
$ws.Range("E12").Value = 1.633
$ws.Range("F12").Value = 1.032
$ws.Range("G12").Value = 1.016
$ws.Range("H12").Value = 1.556
$ws.Range("I12").Value = 1.191
$ws.Range("J12").Value = 1.091
$ws.Range("K12").Value = 1.407
$ws.Range("L12").Value = 0.538
$ws.Range("M12").Value = 0.733
$ws.Range("N12").Value = 1.389
$ws.Range("O12").Value = 0.571
$ws.Range("P12").Value = 0.756
$ws.Range("Z12").Value = 1.25
$ws.Range("AA13").Value = 3.901
$ws.Range("AB13").Value = 1.975
$ws.Range("AC13").Value = 6.273
$ws.Range("AD13").Value = 2.88
$ws.Range("AE13").Value = 1.697
$ws.Range("AF13").Value = 1.667
$ws.Range("AG13").Value = 0.8
$ws.Range("AH13").Value = 0.894
$ws.Range("AI13").Value = 1.311
$ws.Range("AJ13").Value = 0.348
$ws.Range("AK13").Value = 0.59
$ws.Range("AL13").Value = 1.689
$ws.Range("AM13").Value = 0.792
$ws.Range("AN13").Value = 0.89
$ws.Range("AO13").Value = 1.556
$ws.Range("B13").Value = 3.533
$ws.Range("C13").Value = 1.404
$ws.Range("D13").Value = 1.185
$ws.Range("E13").Value = 4.564
$ws.Range("F13").Value = 0.707
$ws.Range("G13").Value = 0.841
$ws.Range("H13").Value = 4.524
$ws.Range("I13").Value = 0.916
$ws.Range("J13").Value = 0.957
$ws.Range("K13").Value = 2.3
$ws.Range("L13").Value = 0.61
$ws.Range("M13").Value = 0.781
$ws.Range("N13").Value = 3.333
$ws.Range("O13").Value = 0.756
$ws.Range("P13").Value = 0.869
$ws.Range("Z13").Value = 2.833
